# Add a new "email: 100" field paragraph right after the existing
# "passwordHash: 100" field paragraph (mirrors the formatting/structure
# used by the other field paragraphs: a leading tab followed by the
# "<name>: <length>" text, with de-DE language formatting).

$d = $word.ActiveDocument

# Locate the paragraph that holds "passwordHash: 100" so the new
# paragraph can be inserted immediately after it, right before the
# "phoneNumberHash: 100" paragraph.
$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -match "passwordHash") {
        $anchorIndex = $i
        break
    }
}

$anchorRange = $d.Paragraphs.Item($anchorIndex).Range
$anchorRange.Collapse(0)
$anchorRange.InsertParagraphAfter()

# The freshly inserted (still empty) paragraph is now the next one.
$newRange = $d.Paragraphs.Item($anchorIndex + 1).Range

# Build the new run as real OOXML (tab element + text run) so it matches
# the same <w:tab/><w:t>...</w:t> pattern used elsewhere in the document,
# rather than a literal tab character inside the text run.
$newParaXml = '<w:p>' +
    '<w:pPr><w:rPr><w:lang w:val="de-DE"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:tab/><w:t>email: 100</w:t></w:r>' +
  '</w:p>'

$packageXml = '<?xml version="1.0" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
      '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
          '<w:body>' + $newParaXml + '</w:body>' +
        '</w:document>' +
      '</pkg:xmlData>' +
    '</pkg:part>' +
  '</pkg:package>'

[void]$newRange.InsertXML($packageXml)

Write-Output "Inserted 'email: 100' paragraph after paragraph $anchorIndex."
